# Applies the recalculated "Leve Profit" tracker values produced by the
# scheduled runner (chore: update Sheets via scheduled runner).
# For each affected row, the currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ /
# LeveProfitHQ columns (H:N) are refreshed with newer market-board figures.
# A few rows also gain or lose trailing columns entirely (not just values),
# which is handled via ClearContents() for removed cells and a plain
# Value assignment for newly populated cells.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 14505.75
$ws.Range("I18").Value = 17507.334
$ws.Range("K18").Value = 17507.334
$ws.Range("M18").Value = -17223.334
$ws.Range("H40").Value = 2205
$ws.Range("J40").Value = 2279.4
$ws.Range("L40").Value = 2279.4
$ws.Range("N40").Value = -2629.4
$ws.Range("H55").Value = 393.2857
$ws.Range("I55").Value = 518
$ws.Range("J55").Value = 299.75
$ws.Range("K55").Value = 518
$ws.Range("L55").Value = 299.75
$ws.Range("M55").Value = -304
$ws.Range("N55").Value = -727.75
$ws.Range("H58").Value = 2148
$ws.Range("I58").Value = 947.3333
$ws.Range("J58").Value = 5750
$ws.Range("K58").Value = 2841.9999
$ws.Range("L58").Value = 17250
$ws.Range("M58").Value = -2691.9999
$ws.Range("N58").Value = -17550
$ws.Range("H92").Value = 914.2
$ws.Range("I92").Value = 704.625
$ws.Range("J92").Value = 1752.5
$ws.Range("K92").Value = 704.625
$ws.Range("L92").Value = 1752.5
$ws.Range("M92").Value = 543.375
$ws.Range("N92").Value = -4248.5
$ws.Range("H100").Value = 1956.7858
$ws.Range("I100").Value = 1449.6666
$ws.Range("K100").Value = 1449.6666
$ws.Range("M100").Value = -908.6666
$ws.Range("H116").Value = 4073.0557
$ws.Range("I116").Value = 3419.182
$ws.Range("J116").Value = 5100.5713
$ws.Range("K116").Value = 3419.182
$ws.Range("L116").Value = 5100.5713
$ws.Range("M116").Value = 22.81800000000021
$ws.Range("N116").Value = -11984.5713
$ws.Range("H137").Value = 4945.857
$ws.Range("I137").Value = 4353.5
$ws.Range("K137").Value = 13060.5
$ws.Range("M137").Value = -10510.5
$ws.Range("H138").Value = 6572.9775
$ws.Range("I138").Value = 3204.5715
$ws.Range("J138").Value = 7193.4736
$ws.Range("K138").Value = 9613.7145
$ws.Range("L138").Value = 21580.4208
$ws.Range("M138").Value = -4473.7145
$ws.Range("N138").Value = -31860.4208

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17188
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("H132").Value = 17375.25
$ws.Range("I132").Value = 21373.75
$ws.Range("J132").Value = 9378.25
$ws.Range("K132").Value = 64121.25
$ws.Range("L132").Value = 28134.75
$ws.Range("M132").Value = -61591.25
$ws.Range("N132").Value = -33194.75
$ws.Range("N32").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 258.6
$ws.Range("I19").Value = 123.25
$ws.Range("K19").Value = 123.25
$ws.Range("M19").Value = 46.75
$ws.Range("H24").Value = 258.6
$ws.Range("I24").Value = 123.25
$ws.Range("K24").Value = 123.25
$ws.Range("M24").Value = 46.75
$ws.Range("H31").Value = 11699.55
$ws.Range("I31").Value = 15199.5
$ws.Range("J31").Value = 8199.6
$ws.Range("K31").Value = 15199.5
$ws.Range("L31").Value = 8199.6
$ws.Range("M31").Value = -14904.5
$ws.Range("N31").Value = -8789.6
$ws.Range("H34").Value = 11699.55
$ws.Range("I34").Value = 15199.5
$ws.Range("J34").Value = 8199.6
$ws.Range("K34").Value = 15199.5
$ws.Range("L34").Value = 8199.6
$ws.Range("M34").Value = -14997.5
$ws.Range("N34").Value = -8603.6
$ws.Range("H41").Value = 15000
$ws.Range("I41").Value = 15000
$ws.Range("K41").Value = 15000
$ws.Range("M41").Value = -14572
$ws.Range("H60").Value = 10792.923
$ws.Range("I60").Value = 8775.666999999999
$ws.Range("J60").Value = 35000
$ws.Range("K60").Value = 8775.666999999999
$ws.Range("L60").Value = 35000
$ws.Range("M60").Value = -8264.666999999999
$ws.Range("H62").Value = 3526.25
$ws.Range("I62").Value = 2052.5
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 2052.5
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -1428.5
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 3526.25
$ws.Range("I65").Value = 2052.5
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 10262.5
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -7142.5
$ws.Range("N65").Value = -31240
$ws.Range("H74").Value = 37999
$ws.Range("I74").Value = 9498.5
$ws.Range("K74").Value = 9498.5
$ws.Range("M74").Value = -8624.5
$ws.Range("H77").Value = 37999
$ws.Range("I77").Value = 9498.5
$ws.Range("K77").Value = 28495.5
$ws.Range("M77").Value = -24127.5
$ws.Range("N60").Value = -36022

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 170108.23
$ws.Range("I2").Value = 137530.75
$ws.Range("J2").Value = 222232.2
$ws.Range("K2").Value = 825184.5
$ws.Range("L2").Value = 1333393.2
$ws.Range("M2").Value = -825071.5
$ws.Range("N2").Value = -1333619.2
$ws.Range("H46").Value = 2114.1428
$ws.Range("I46").Value = 1066.3334
$ws.Range("K46").Value = 3199.0002
$ws.Range("M46").Value = -3108.0002
$ws.Range("H63").Value = 14320.333
$ws.Range("I63").Value = 14320.333
$ws.Range("K63").Value = 42960.999
$ws.Range("M63").Value = -42211.999
$ws.Range("H66").Value = 14320.333
$ws.Range("I66").Value = 14320.333
$ws.Range("K66").Value = 128882.997
$ws.Range("M66").Value = -125138.997
$ws.Range("H97").Value = 5210312.5
$ws.Range("I97").Value = 994
$ws.Range("J97").Value = 8931254
$ws.Range("K97").Value = 2982
$ws.Range("L97").Value = 26793762
$ws.Range("M97").Value = -2486
$ws.Range("N97").Value = -26794754
$ws.Range("H113").Value = 2635.15
$ws.Range("I113").Value = 2596.3333
$ws.Range("K113").Value = 7788.999899999999
$ws.Range("M113").Value = -5618.999899999999
$ws.Range("H141").Value = 11999.5
$ws.Range("I141").Value = 11999.5
$ws.Range("K141").Value = 35998.5
$ws.Range("M141").Value = -30818.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 103399.8
$ws.Range("J125").Value = 103399.8
$ws.Range("L125").Value = 103399.8
$ws.Range("N125").Value = -108319.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1829.0952
$ws.Range("I22").Value = 1730.7646
$ws.Range("J22").Value = 2247
$ws.Range("K22").Value = 1730.7646
$ws.Range("L22").Value = 2247
$ws.Range("M22").Value = -1435.7646
$ws.Range("N22").Value = -2837
$ws.Range("H27").Value = 1829.0952
$ws.Range("I27").Value = 1730.7646
$ws.Range("J27").Value = 2247
$ws.Range("K27").Value = 1730.7646
$ws.Range("L27").Value = 2247
$ws.Range("M27").Value = -1623.7646
$ws.Range("N27").Value = -2461
$ws.Range("H46").Value = 2796.3635
$ws.Range("I46").Value = 2751.4285
$ws.Range("J46").Value = 2875
$ws.Range("K46").Value = 2751.4285
$ws.Range("L46").Value = 2875
$ws.Range("M46").Value = -2563.4285
$ws.Range("N46").Value = -3251
$ws.Range("H68").Value = 2239.2
$ws.Range("J68").Value = 1999
$ws.Range("L68").Value = 1999
$ws.Range("N68").Value = -3497
$ws.Range("H71").Value = 2239.2
$ws.Range("J71").Value = 1999
$ws.Range("L71").Value = 9995
$ws.Range("N71").Value = -17483

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 2436.3333
$ws.Range("I26").Value = 923.6
$ws.Range("K26").Value = 923.6
$ws.Range("M26").Value = -630.6
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 3500
$ws.Range("K132").Value = 10500
$ws.Range("M132").Value = -7970
$ws.Range("H136").Value = 3581.3125
$ws.Range("J136").Value = 3852
$ws.Range("L136").Value = 11556
$ws.Range("N136").Value = -16656
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()
